$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a new "2021" column (R) that mirrors the formatting of the existing
# "2020" column (Q), then overwrite the values for the new year.
$ws.Range("Q4:Q14").Copy($ws.Range("R4:R14"))

$ws.Range("R4").Value = 2021
$ws.Range("R5").Value = 5.8
$ws.Range("R6").Value = 4.7
$ws.Range("R7").Value = 1.6
$ws.Range("R8").Value = 12.9
$ws.Range("R9").Value = 10.199999999999999
$ws.Range("R10").Value = 4.2
$ws.Range("R11").Value = 3.3
$ws.Range("R12").Value = 15.2
$ws.Range("R13").Value = 2.4
$ws.Range("R14").Value = 0.6

# Update the sheet's active selection, matching the recorded edit.
$ws.Range("T9").Select()
